$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = -0.3216975593808727
$ws.Range("J3").Value = 0.6701918066130262
$ws.Range("K3").Value = 0.5263752878107336
$ws.Range("L3").Value = 2.367711971973246

$ws.Range("I20").Value = -0.06444003869336853
$ws.Range("J20").Value = 0.6912239929730917
$ws.Range("K20").Value = 0.1355354538661993
$ws.Range("L20").Value = 2.091729358658431
